$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at the top of the sheet; this pushes the existing
# header row (row 1) and all data rows down by two.
$ws.Rows("1:2").Insert()

# Row 1: report title, merged across A1:F1
$titleCell = $ws.Range("A1")
$titleCell.Value = "Reporte totalidad licencias"
$ws.Range("A1:F1").Merge()
$titleFont = $titleCell.Font
$titleFont.Bold = $true
$titleFont.Size = 24
$titleFont.Color = 16711680
$titleCell.HorizontalAlignment = -4108
$titleCell.VerticalAlignment = -4108

# Row 2: creation date/time, merged across A2:F2
$dateCell = $ws.Range("A2")
$dateCell.Value = "Fecha y hora de creación: 06/07/2021, 22:07:52"
$ws.Range("A2:F2").Merge()
$dateFont = $dateCell.Font
$dateFont.Bold = $true
$dateFont.Size = 14
$dateCell.HorizontalAlignment = -4108
$dateCell.VerticalAlignment = -4108
